$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "70.159.39"
$ws.Range("E2").Value = "  +1.95%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.937.56"
$ws.Range("E3").Value = "  +2.27%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
Set-TextValue "D5" "610.71"
$ws.Range("E5").Value = "  +1.59%  "

# Row 6 - Solana
Set-TextValue "D6" "171.36"
$ws.Range("E6").Value = "  +5.57%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.934.72"
$ws.Range("E7").Value = "  +2.24%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.15%  "

# Row 9 - XRP
Set-TextValue "D9" "0.538"
$ws.Range("E9").Value = "  +1.58%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.75%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +2.26%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.471"
$ws.Range("E12").Value = "  +2.61%  "

# Row 13 - Avalanche
Set-TextValue "D13" "38.66"
$ws.Range("E13").Value = "  +5.05%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +5.78%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.595.85"
$ws.Range("E15").Value = "  +2.28%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.927.98"
$ws.Range("E16").Value = "  +2.24%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "70.152.37"
$ws.Range("E17").Value = "  +1.67%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  +2.03%  "

# Row 19 - Chainlink
Set-TextValue "D19" "18.59"
$ws.Range("E19").Value = "  +8.44%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  -0.90%  "

# Row 21 - Uniswap
Set-TextValue "D21" "11.14"
$ws.Range("E21").Value = "  -2.84%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "497.37"
$ws.Range("E22").Value = "  +2.74%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.749"
$ws.Range("E23").Value = "  +4.15%  "

# Row 24 - PEPE
Set-TextValue "D24" "0.0000167"
$ws.Range("E24").Value = "  +4.94%  "

# Row 25 - Litecoin
Set-TextValue "D25" "86.16"
$ws.Range("E25").Value = "  +2.80%  "

# Row 26 - Fetch.AI
$ws.Range("E26").Value = "  +2.39%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "12.41"
$ws.Range("E27").Value = "  +2.57%  "

# Row 28 - RenderToken
Set-TextValue "D28" "10.21"
$ws.Range("E28").Value = "  +2.11%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.10%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +1.81%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  +3.69%  "

# Row 32 - WrappedeETH
$ws.Range("D32").Value = "4.086.35"

# Row 34 - EthereumClassic
Set-TextValue "D34" "32.38"
$ws.Range("E34").Value = "  +0.55%  "

# Row 35 - RenzoRestakedETH
$ws.Range("D35").Value = "3.898.38"
$ws.Range("E35").Value = "  +2.61%  "

# Row 36 - Hedera
$ws.Range("E36").Value = "  +1.21%  "

# Row 37 - Filecoin
$ws.Range("E37").Value = "  +4.86%  "

# Row 38 - Mantle
$ws.Range("E38").Value = "  +1.54%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +1.64%  "

# Row 41 - TheGraph
Set-TextValue "D41" "0.330"
$ws.Range("E41").Value = "  +3.76%  "

# Row 42 - FirstDigitalUSD
Set-TextValue "D42" "0.999"
$ws.Range("E42").Value = "  -0.05%  "

# Row 43 - Stacks
Set-TextValue "D43" "2.12"
$ws.Range("E43").Value = "  +7.40%  "

# Row 44 - Bittensor
Set-TextValue "D44" "438.56"

# Row 45 - OKB
Set-TextValue "D45" "48.34"
$ws.Range("E45").Value = "  -0.29%  "

# Row 46 - Cosmos
Set-TextValue "D46" "8.69"
$ws.Range("E46").Value = "  +3.77%  "

# Row 48 - VeChain
Set-TextValue "D48" "0.0368"
$ws.Range("E48").Value = "  +3.02%  "

# Row 49 - FLOKI
$ws.Range("E49").Value = "  +22.24%  "

# Row 50 - Arweave
Set-TextValue "D50" "40.80"
$ws.Range("E50").Value = "  +5.72%  "

# Row 51 - Monero
Set-TextValue "D51" "143.54"
$ws.Range("E51").Value = "  +0.34%  "
